$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

# Row 5: add Model annotation in column E
$ws.Range("E5").Value = "Model"

# Row 7 (new): PlacedKingdomItemList field referencing a list-of-packet model
$ws.Range("A7").Value = "PlacedKingdomItemList"
$ws.Range("B7").Value = "LIST:PlacedKingdomItemPacket"
$ws.Range("E7").Value = "Packet"

# Row 4: rename Ysize -> SizeY
$ws.Range("A4").Value = "SizeY"
$ws.Range("B4").Value = "INT"

# Row 3: rename XSize -> SizeX
$ws.Range("A3").Value = "SizeX"
$ws.Range("B3").Value = "INT"

$ws.Range("A4").Select()
